$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "42.798.80"
Set-TextValue $ws.Range("E2") "  -0.35%  "
Set-TextValue $ws.Range("D3") "2.274.97"
Set-TextValue $ws.Range("E3") "  -0.43%  "
Set-TextValue $ws.Range("E4") "  -0.10%  "
Set-TextValue $ws.Range("D5") "250.30"
Set-TextValue $ws.Range("E5") "  -0.56%  "
Set-TextValue $ws.Range("D6") "0.634"
Set-TextValue $ws.Range("E6") "  -1.07%  "
Set-TextValue $ws.Range("D7") "79.27"
Set-TextValue $ws.Range("E7") "  +8.99%  "
Set-TextValue $ws.Range("E8") "  -0.02%  "
Set-TextValue $ws.Range("D9") "0.644"
Set-TextValue $ws.Range("E9") "  -2.53%  "
Set-TextValue $ws.Range("D10") "41.39"
Set-TextValue $ws.Range("E10") "  +5.51%  "
Set-TextValue $ws.Range("D11") "0.0969"
Set-TextValue $ws.Range("E11") "  -1.18%  "
Set-TextValue $ws.Range("D12") "7.39"
Set-TextValue $ws.Range("E12") "  -0.62%  "
Set-TextValue $ws.Range("E13") "  -1.18%  "
Set-TextValue $ws.Range("D14") "2.615.76"
Set-TextValue $ws.Range("E14") "  -0.45%  "
Set-TextValue $ws.Range("D15") "15.13"
Set-TextValue $ws.Range("E15") "  +0.20%  "
Set-TextValue $ws.Range("D16") "0.868"
Set-TextValue $ws.Range("E16") "  -2.99%  "
Set-TextValue $ws.Range("D17") "2.268.54"
Set-TextValue $ws.Range("E17") "  -1.61%  "
Set-TextValue $ws.Range("D18") "42.708.15"
Set-TextValue $ws.Range("E18") "  -0.39%  "
Set-TextValue $ws.Range("E19") "  -1.46%  "
Set-TextValue $ws.Range("D20") "6.23"
Set-TextValue $ws.Range("E20") "  -2.58%  "
Set-TextValue $ws.Range("D21") "72.11"
Set-TextValue $ws.Range("E21") "  -2.07%  "
Set-TextValue $ws.Range("D22") "233.75"
Set-TextValue $ws.Range("E22") "  -1.71%  "
Set-TextValue $ws.Range("E23") "  -0.15%  "
Set-TextValue $ws.Range("D24") "3.79"
Set-TextValue $ws.Range("E24") "  -2.54%  "
Set-TextValue $ws.Range("D26") "11.32"
Set-TextValue $ws.Range("E26") "  -4.29%  "
Set-TextValue $ws.Range("E27") "  -4.66%  "
Set-TextValue $ws.Range("E28") "  +1.80%  "
Set-TextValue $ws.Range("D29") "169.99"
Set-TextValue $ws.Range("E29") "  +1.08%  "
Set-TextValue $ws.Range("D30") "20.86"
Set-TextValue $ws.Range("E31") "  +5.14%  "
Set-TextValue $ws.Range("D32") "0.0849"
Set-TextValue $ws.Range("E32") "  +4.27%  "
Set-TextValue $ws.Range("D33") "0.122"
Set-TextValue $ws.Range("E33") "  -4.85%  "
Set-TextValue $ws.Range("D34") "30.68"
Set-TextValue $ws.Range("E34") "  -2.46%  "
Set-TextValue $ws.Range("E35") "  +0.40%  "
Set-TextValue $ws.Range("D36") "4.56"
Set-TextValue $ws.Range("E36") "  -5.63%  "
Set-TextValue $ws.Range("D37") "4.77"
Set-TextValue $ws.Range("E37") "  -0.62%  "
Set-TextValue $ws.Range("D38") "0.0303"
Set-TextValue $ws.Range("D39") "13.52"
Set-TextValue $ws.Range("E39") "  +1.12%  "
Set-TextValue $ws.Range("D40") "2.26"
Set-TextValue $ws.Range("E40") "  -2.85%  "
Set-TextValue $ws.Range("D41") "5.94"
Set-TextValue $ws.Range("E41") "  -2.21%  "
Set-TextValue $ws.Range("D42") "115.34"
Set-TextValue $ws.Range("E42") "  +18.11%  "
Set-TextValue $ws.Range("D43") "0.209"
Set-TextValue $ws.Range("E43") "  -1.56%  "
Set-TextValue $ws.Range("D44") "61.23"
Set-TextValue $ws.Range("E44") "  -1.16%  "
Set-TextValue $ws.Range("E45") "  -3.26%  "
Set-TextValue $ws.Range("E46") "  -2.02%  "
Set-TextValue $ws.Range("D47") "4.64"
Set-TextValue $ws.Range("E47") "  -6.53%  "
Set-TextValue $ws.Range("E48") "  -0.15%  "
Set-TextValue $ws.Range("E49") "  -3.96%  "
Set-TextValue $ws.Range("D50") "1.17"
Set-TextValue $ws.Range("E50") "  -2.32%  "
Set-TextValue $ws.Range("D51") "4.27"
Set-TextValue $ws.Range("E51") "  -1.53%  "
